# Fix duplicate name issue and update gender identification.
# - Standardize 'Solum Ole Peder Uthus' -> remove duplicate row on Male_50m,
#   shifting the following row up and adding the new "Kristian Volden" entry.
# - Move 'Bjarne Forfot' (misclassified) out of Female_50m, replacing the
#   row with the correct next female entry "Miriam Vedvik".

$wb = $excel.ActiveWorkbook

# --- Male_50m sheet ---
$maleWs = $wb.Worksheets.Item("Male_50m")

# Dato column holds plain text like "10.06.2018" - force text format so
# Excel doesn't reinterpret the strings as date serials.
$maleWs.Range("D10:D11").NumberFormat = "@"

# Row 10 becomes what used to be row 11 (Odin Spangen Normann), removing the
# duplicate "Solum Ole Peder Uthus" entry.
$maleWs.Cells.Item(10, 1).Value = "Odin Spangen Normann"
$maleWs.Cells.Item(10, 2).Value = "2.26,78"
$maleWs.Cells.Item(10, 3).Value = 424
$maleWs.Cells.Item(10, 4).Value = "04.03.2012"
$maleWs.Cells.Item(10, 5).Value = "Berlin"

# Row 11 becomes the new entry, "Kristian Volden".
$maleWs.Cells.Item(11, 1).Value = "Kristian Volden"
$maleWs.Cells.Item(11, 2).Value = "2.28,14"
$maleWs.Cells.Item(11, 3).Value = 413
$maleWs.Cells.Item(11, 4).Value = "28.05.2016"
$maleWs.Cells.Item(11, 5).Value = "Bergen"

# --- Female_50m sheet ---
$femaleWs = $wb.Worksheets.Item("Female_50m")

$femaleWs.Range("D11").NumberFormat = "@"

# Row 11 "Bjarne Forfot" (male, misclassified) is replaced by the next
# correct female entry "Miriam Vedvik".
$femaleWs.Cells.Item(11, 1).Value = "Miriam Vedvik"
$femaleWs.Cells.Item(11, 2).Value = "3.12,55"
$femaleWs.Cells.Item(11, 3).Value = 253
$femaleWs.Cells.Item(11, 4).Value = "28.04.2007"
$femaleWs.Cells.Item(11, 5).Value = "Namsos"
